# update read file excel in appyTestCase.java
#
# - row 8 on "data": the blank G8 cell is filled in with the same
#   "This is the test content" comment used by the other data rows, and
#   the stray H8 phone-number value that row 8 was still carrying is
#   cleared out so the row lines up with its siblings.
# - the active selection on "data" moves from D2 to L4.
# - a brand-new, empty worksheet named "Element" is appended after "data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 8: fill G8 with the shared "This is the test content" string and
# drop the leftover H8 number.
$ws.Range("G8").Value = "This is the test content"
$ws.Range("H8").ClearContents()

# Move the selection on the "data" sheet to L4.
$ws.Range("L4").Select()

# Append a new, empty worksheet named "Element" right after "data".
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$newSheet.Name = "Element"

# Leave "data" as the active/selected sheet (matches the unchanged
# tabSelected="1" on the "data" sheetView in the target workbook).
$ws.Activate()
$ws.Select()
